$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B width (closest achievable to target 26.140625 given runtime width granularity)
$ws.Columns.Item(2).ColumnWidth = 25.25

# Update column D (Label2) values for rows with new label data
$d = @{
    3 = 4
    4 = 4
    5 = 1
    6 = 3
    7 = 4
    12 = 4
    14 = 3
    15 = 2
    17 = 2
    18 = 1
    19 = 4
    21 = 3
    22 = 2
    23 = 4
    24 = 3
    25 = 2
    26 = 3
    31 = 1
    32 = 2
    34 = 1
    35 = 3
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    50 = 3
    51 = 1
    52 = 1
    53 = 3
    54 = 3
    55 = 3
    56 = 3
    57 = 1
    58 = 1
    59 = 1
    60 = 1
    61 = 1
    62 = 1
    63 = 1
    64 = 3
    65 = 3
    66 = 3
    68 = 1
    69 = 3
    70 = 1
    71 = 3
    72 = 3
    74 = 1
    75 = 3
    76 = 3
    77 = 1
    78 = 3
    79 = 1
    80 = 1
    82 = 1
    87 = 1
    89 = 1
    90 = 1
    91 = 1
    92 = 1
    93 = 1
    95 = 1
    96 = 2
    97 = 2
    99 = 3
    100 = 3
    101 = 3
    102 = 3
    103 = 3
    104 = 3
    105 = 3
    106 = 3
    107 = 1
    108 = 3
    109 = 3
    110 = 3
    111 = 3
    112 = 2
    113 = 2
    114 = 2
    115 = 3
    116 = 2
    118 = 1
    119 = 1
    120 = 1
    124 = 1
    125 = 1
    126 = 1
    127 = 1
    128 = 1
    129 = 1
    130 = 1
    131 = 1
    135 = 1
    138 = 1
    140 = 1
    141 = 1
    144 = 1
    145 = 1
    147 = 3
    148 = 1
    149 = 1
    150 = 2
    151 = 1
    152 = 2
    153 = 1
    154 = 1
    155 = 2
    156 = 1
    157 = 1
    162 = 3
    163 = 1
    164 = 1
    165 = 3
    166 = 1
    168 = 3
    169 = 4
    171 = 1
    174 = 1
}

foreach ($row in $d.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $d[$row]
}

# Move active selection to D2, matching the final sheet view state
$ws.Range("D2").Select()
